# Break out Equipment and Location changes
#
# The "IMC" (equipment) and "Stock Number" (location) values recorded on
# Sheet1 are updated from the old combined code "A123" to the new code
# "W333" for both data rows. Updating the cell values causes the shared
# strings table to drop the now-unused "A123"/"A123 200017758" entries and
# append the new "W333"/"W333 200017758" entries, matching the recorded
# change. The active selection is also moved from the old row-4 range
# selection down to a single cell, I4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "W333"
$ws.Range("I2").Value = "W333 200017758"

$ws.Range("H3").Value = "W333"
$ws.Range("I3").Value = "W333 200017758"

$ws.Range("I4").Select()
